# Scheduled data refresh: update cached market-price / profit figures
# on the Leve profit-tracking sheets (columns H-N are pasted values, not formulas).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 233657.39
$ws.Range("J43").Value = 587638.4399999999
$ws.Range("L43").Value = 587638.4399999999
$ws.Range("N43").Value = -587776.4399999999
# Row 70
$ws.Range("H70").Value = 97225460
$ws.Range("I70").Value = 83334130
$ws.Range("K70").Value = 250002390
$ws.Range("M70").Value = -250002120
# Row 73
$ws.Range("H73").Value = 97225460
$ws.Range("I73").Value = 83334130
$ws.Range("K73").Value = 250002390
$ws.Range("M73").Value = -250001454
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").ClearContents()
$ws.Range("N87").Value = 0
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").ClearContents()
$ws.Range("N90").Value = 0
# Row 98
$ws.Range("H98").Value = 8715.833000000001
$ws.Range("I98").Value = 9373.75
$ws.Range("J98").Value = 7400
$ws.Range("K98").Value = 9373.75
$ws.Range("L98").Value = 7400
$ws.Range("M98").Value = -7875.75
$ws.Range("N98").Value = -10396
# Row 103
$ws.Range("H103").Value = 3179.2
$ws.Range("I103").Value = 2000
$ws.Range("K103").Value = 6000
$ws.Range("M103").Value = -5414
# Row 112
$ws.Range("H112").Value = 10349.695
$ws.Range("J112").Value = 10729.228
$ws.Range("L112").Value = 32187.684
$ws.Range("N112").Value = -34403.68399999999
# Row 122
$ws.Range("H122").Value = 8715.833000000001
$ws.Range("I122").Value = 9373.75
$ws.Range("J122").Value = 7400
$ws.Range("K122").Value = 28121.25
$ws.Range("L122").Value = 22200
$ws.Range("M122").Value = -25671.25
$ws.Range("N122").Value = -27100
# Row 132
$ws.Range("H132").Value = 1848.2222
$ws.Range("I132").Value = 1805.8
$ws.Range("J132").Value = 3333
$ws.Range("K132").Value = 5417.4
$ws.Range("L132").Value = 9999
$ws.Range("M132").Value = -2887.4
$ws.Range("N132").Value = -15059
# Row 137
$ws.Range("H137").Value = 3314.3333
$ws.Range("I137").Value = 4245.923
$ws.Range("J137").Value = 2449.2856
$ws.Range("K137").Value = 12737.769
$ws.Range("L137").Value = 7347.8568
$ws.Range("M137").Value = -10187.769
$ws.Range("N137").Value = -12447.8568

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4154.6294
$ws.Range("I32").Value = 4455.841
$ws.Range("K32").Value = 4455.841
$ws.Range("M32").Value = -4168.841
# Row 61
$ws.Range("H61").Value = 5133.3438
$ws.Range("I61").Value = 2791.9412
$ws.Range("J61").Value = 14318.846
$ws.Range("K61").Value = 2791.9412
$ws.Range("L61").Value = 14318.846
$ws.Range("M61").Value = -2579.9412
$ws.Range("N61").Value = -14742.846
# Row 122
$ws.Range("H122").Value = 10433.689
$ws.Range("I122").Value = 18903.846
$ws.Range("K122").Value = 56711.538
$ws.Range("M122").Value = -54261.538
# Row 132
$ws.Range("H132").Value = 1388292.2
$ws.Range("I132").Value = 1646125.2
$ws.Range("J132").Value = 13183.333
$ws.Range("K132").Value = 4938375.6
$ws.Range("L132").Value = 39549.999
$ws.Range("M132").Value = -4935845.6
$ws.Range("N132").Value = -44609.999
# Row 136
$ws.Range("H136").Value = 5133.3438
$ws.Range("I136").Value = 2791.9412
$ws.Range("J136").Value = 14318.846
$ws.Range("K136").Value = 8375.8236
$ws.Range("L136").Value = 42956.538
$ws.Range("M136").Value = -5825.8236
$ws.Range("N136").Value = -48056.538

$ws = $wb.Worksheets.Item("BSM")
# Row 50
$ws.Range("H50").Value = 45435.332
$ws.Range("I50").Value = 38460
$ws.Range("J50").Value = 48923
$ws.Range("K50").Value = 38460
$ws.Range("L50").Value = 48923
$ws.Range("M50").Value = -37886
$ws.Range("N50").Value = -50071
# Row 105
$ws.Range("H105").Value = 3647.75
$ws.Range("I105").Value = 1261.8334
$ws.Range("K105").Value = 1261.8334
$ws.Range("M105").Value = 485.1666
# Row 134
$ws.Range("H134").Value = 5501.643
$ws.Range("I134").Value = 1170.375
$ws.Range("J134").Value = 11276.667
$ws.Range("K134").Value = 3511.125
$ws.Range("L134").Value = 33830.001
$ws.Range("M134").Value = -976.125
$ws.Range("N134").Value = -38900.001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6406.089
$ws.Range("I31").Value = 3161.6
$ws.Range("J31").Value = 10461.7
$ws.Range("K31").Value = 3161.6
$ws.Range("L31").Value = 10461.7
$ws.Range("M31").Value = -2866.6
$ws.Range("N31").Value = -11051.7
# Row 34
$ws.Range("H34").Value = 6406.089
$ws.Range("I34").Value = 3161.6
$ws.Range("J34").Value = 10461.7
$ws.Range("K34").Value = 3161.6
$ws.Range("L34").Value = 10461.7
$ws.Range("M34").Value = -2959.6
$ws.Range("N34").Value = -10865.7
# Row 39
$ws.Range("H39").Value = 6999
$ws.Range("I39").Value = 6999
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 6999
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -6608
# Row 49
$ws.Range("H49").Value = 6999
$ws.Range("I49").Value = 6999
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 6999
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -6817
# Row 99
$ws.Range("H99").Value = 10571.143
$ws.Range("I99").Value = 19333
$ws.Range("J99").Value = 8181.5454
$ws.Range("K99").Value = 19333
$ws.Range("L99").Value = 8181.5454
$ws.Range("M99").Value = -17835
$ws.Range("N99").Value = -11177.5454
# Row 122
$ws.Range("H122").Value = 3017.5
$ws.Range("I122").Value = 3177.1428
$ws.Range("K122").Value = 9531.428400000001
$ws.Range("M122").Value = -7081.428400000001
# Row 126
$ws.Range("H126").Value = 10571.143
$ws.Range("I126").Value = 19333
$ws.Range("J126").Value = 8181.5454
$ws.Range("K126").Value = 57999
$ws.Range("L126").Value = 24544.6362
$ws.Range("M126").Value = -55529
$ws.Range("N126").Value = -29484.6362
# Row 141
$ws.Range("H141").Value = 63352
$ws.Range("J141").Value = 63352
$ws.Range("L141").Value = 63352
$ws.Range("N141").Value = -73712

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 6180
$ws.Range("J34").Value = 8944.546
$ws.Range("L34").Value = 26833.638
$ws.Range("N34").Value = -27001.638
# Row 41
$ws.Range("H41").Value = 1166.6666
$ws.Range("I41").Value = 1000
$ws.Range("J41").Value = 1500
$ws.Range("K41").Value = 3000
$ws.Range("L41").Value = 4500
$ws.Range("M41").Value = -2662
$ws.Range("N41").Value = -5176
# Row 43
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").ClearContents()
$ws.Range("N43").Value = 0
# Row 74
$ws.Range("H74").Value = 2987
$ws.Range("J74").Value = 2987
$ws.Range("L74").Value = 8961
$ws.Range("N74").Value = -11083
# Row 77
$ws.Range("H77").Value = 2987
$ws.Range("J77").Value = 2987
$ws.Range("L77").Value = 26883
$ws.Range("N77").Value = -37491
# Row 104
$ws.Range("H104").Value = 4944.8
$ws.Range("J104").Value = 4944.8
$ws.Range("L104").Value = 14834.4
$ws.Range("N104").Value = -20076.4

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 9211.625
$ws.Range("J70").Value = 12000
$ws.Range("L70").Value = 12000
$ws.Range("N70").Value = -12540
# Row 73
$ws.Range("H73").Value = 9211.625
$ws.Range("J73").Value = 12000
$ws.Range("L73").Value = 12000
$ws.Range("N73").Value = -13872
# Row 97
$ws.Range("H97").Value = 1335.4
$ws.Range("I97").Value = 1157.4445
$ws.Range("J97").Value = 2937
$ws.Range("K97").Value = 1157.4445
$ws.Range("L97").Value = 2937
$ws.Range("M97").Value = -661.4445000000001
$ws.Range("N97").Value = -3929
# Row 122
$ws.Range("H122").Value = 111118200
$ws.Range("I122").Value = 333338020
$ws.Range("K122").Value = 1000014060
$ws.Range("M122").Value = -1000011610
# Row 132
$ws.Range("H132").Value = 4029.8965
$ws.Range("I132").Value = 2187.1924
$ws.Range("K132").Value = 6561.5772
$ws.Range("M132").Value = -4031.5772

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 5010.1
$ws.Range("I40").Value = 4412.3335
$ws.Range("J40").Value = 5906.75
$ws.Range("K40").Value = 4412.3335
$ws.Range("L40").Value = 5906.75
$ws.Range("M40").Value = -4276.3335
$ws.Range("N40").Value = -6178.75

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 23363838
$ws.Range("I81").Value = 2501409.2
$ws.Range("J81").Value = 40053780
$ws.Range("K81").Value = 5002818.4
$ws.Range("L81").Value = 80107560
$ws.Range("M81").Value = -5001757.4
$ws.Range("N81").Value = -80109682
# Row 84
$ws.Range("H84").Value = 23363838
$ws.Range("I84").Value = 2501409.2
$ws.Range("J84").Value = 40053780
$ws.Range("K84").Value = 25014092
$ws.Range("L84").Value = 400537800
$ws.Range("M84").Value = -25008788
$ws.Range("N84").Value = -400548408
# Row 122
$ws.Range("H122").Value = 101021.44
$ws.Range("I122").Value = 150934.86
$ws.Range("K122").Value = 452804.58
$ws.Range("M122").Value = -450354.58
# Row 139
$ws.Range("H139").Value = 89999
$ws.Range("J139").Value = 89999
$ws.Range("L139").Value = 89999
$ws.Range("N139").Value = -100279
# Row 141
$ws.Range("H141").Value = 90000
$ws.Range("J141").Value = 90000
$ws.Range("L141").Value = 90000
